# "Generate Report for Handback" — a handback of e2e/a.md was processed and
# the per-locale status rows are refreshed to reflect the new handback:
#   - Status moves to "Handed back: in sync with en-US" (zh-cn only — the
#     de-de row was already "In Translation" and stays that way)
#   - Latest Handback DateTime / Name are stamped with the new handback
#   - The stale "version not latest" Error Detail is cleared now that the
#     handback is in sync

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet : row 2 is the a.md entry ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("L2").Value = "2017-02-22 07:08:07"
$zh.Range("M2").Value = "TestHandback_201702220307"
$zh.Range("R2").Value = ""

# ---- de-de sheet : row 2 is the a.md entry ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("L2").Value = "2017-02-22 07:08:28"
$de.Range("M2").Value = "TestHandback_201702220307"
$de.Range("R2").Value = ""

# ---- widen the now-longer Status / Handback-DateTime / Handback-Name
#      columns so the new text is readable ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

$zh.Columns.Item(3).ColumnWidth = 29.2
$zh.Columns.Item(13).ColumnWidth = 27.17

$de.Columns.Item(3).ColumnWidth = 29.2
$de.Columns.Item(13).ColumnWidth = 27.17
